$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# usuario: autotest11 -> invictus10
$ws.Range("D2").Value = "invictus10"

# convenio: "65401" (text) -> 65437 (number)
# Copy the number format (General) from N2 so the cell is stored numerically
# instead of as a shared text string, matching style index 3.
$ws.Range("N2").Copy()
$ws.Range("M2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M2").Value = 65437

# descripcionFactura: "descripcion 1" -> "Factura auto"
$ws.Range("N2").Value = "Factura auto"
